$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "github username" column (E) for the rows that previously held
# the placeholder "???" value. Cells are written in the same order the
# author appears to have typed them (row 21, 15, 12, 20, 6, 19, 4) so that
# new shared-string entries are appended in that order.
$ws.Range("E21").Value = "BENKECHIKECH"
$ws.Range("E15").Value = "Aminechakr"
$ws.Range("E12").Value = "MedEIP"
$ws.Range("E20").Value = "IssamAxaTech "
$ws.Range("E6").Value  = "braadil"
$ws.Range("E19").Value = "MeriemHamdaoui"
$ws.Range("E4").Value  = "hrsanaa"

# Column E got visibly wider once "MeriemHamdaoui" (the longest new value)
# was entered - mirror Excel's best-fit recalculation.
$ws.Columns.Item(5).ColumnWidth = 19.86

# Re-point the autofilter (and its backing defined name) to cover the
# whole table instead of just the header row.
$ws.AutoFilterMode = $false
[void]$ws.Range("B3:E21").AutoFilter()

$name = $wb.Names.Item(1)
$name.RefersTo = "=Feuil1!`$B`$3:`$E`$21"

# Drop the stale cell selection left over from editing E15.
[void]$ws.Range("A1").Select()
